# Add data for 2021-10-27 (updates "through 10-18" snapshot to "through 10-19")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-10-19"

# Update the October row label
$ws.Range("A12").Value = "October (through 10-19)"

# Row 12 (October partial-month row) updates
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 0.0556
$ws.Range("F12").Value = 27
$ws.Range("G12").Value = 0.1
$ws.Range("I12").Value = 25
$ws.Range("J12").Value = 0.2188
$ws.Range("L12").Value = 44
$ws.Range("M12").Value = 0.0638
$ws.Range("R12").Value = 90
$ws.Range("U12").Value = 124

# Row 13 (Total) updates
$ws.Range("C13").Value = 213
$ws.Range("D13").Value = 0.127
$ws.Range("F13").Value = 410
$ws.Range("G13").Value = 0.1068
$ws.Range("I13").Value = 602
$ws.Range("J13").Value = 0.0865
$ws.Range("L13").Value = 531
$ws.Range("M13").Value = 0.1076
$ws.Range("R13").Value = 938
$ws.Range("S13").Value = 0.0535
$ws.Range("U13").Value = 1289
$ws.Range("V13").Value = 0.0605
